# Commit: "Tue, Jul 07, 2020  2:05:55 PM"
#
# The deck currently has its live/applied design theme (the theme used by the
# slide master and referenced at the presentation level) set to the
# "Integral" color scheme, while an older, unused "Office Theme" color
# scheme is only kept around (attached to the Notes Master). The edit swaps
# the presentation's applied theme colors from the "Integral" palette over
# to the standard "Office" palette -- i.e. the design's theme color scheme
# is changed from Integral to Office.
#
# PowerPoint's Theme Colors are exposed as a 12-slot indexed collection
# (Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink) via
# ThemeColorScheme on the slide master's Theme. We reassign every slot's RGB
# value to the corresponding "Office" theme color.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# msoThemeDark1 / msoThemeLight1 / msoThemeDark2 / msoThemeLight2
$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6

# msoThemeAccent1 .. msoThemeAccent6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47

# msoThemeHyperlink / msoThemeFollowedHyperlink
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
